$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing transcript rows (B2:B6) with the corrected/expanded text.
$ws.Cells.Item(2, 2).Value = "  А можете выделить какие у этих фигур есть признаки?  По которым они варьируются?  Признаки? Размер и метка."
$ws.Cells.Item(3, 2).Value = "  Размер и метка.  А еще какие-нибудь есть?  Нет.  Ну, это...  Они одинаковые по названию одной игрушки."
$ws.Cells.Item(4, 2).Value = "  Одинаковый размер и метка.  А в принципе, ну, вот рандомная фигура, какие можно признаки выделить, если ее не относили к группе?"
$ws.Cells.Item(5, 2).Value = "  А, чтобы в новую группу определить размер и присутствие или отсутствие метки.  А если не определять? Просто как абстрактная фигура?  Любая абстрактная фигура."
$ws.Cells.Item(6, 2).Value = "  Тут имеет значение только метка и размер фигуры, то есть масштаб признаков.  То есть тут есть маленькие фигуры."

# Extend the table with rows 7-14 (A: running index 5-12 with the same style
# as the existing index column, B: the corresponding transcript chunk).
$ws.Range("A6").Copy()
$ws.Range("A7:A14").PasteSpecial(-4122)

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(14, 1).Value = 12

$ws.Cells.Item(7, 2).Value = "  Большие фигуры с меткой, большие фигуры без метки, большие фигуры с меткой и средние, хотя нет, это тоже, наверное, маленькие фигуры без метки.  Так..."
$ws.Cells.Item(8, 2).Value = "  Дополнительные фигуры не вызывали сомнения, куда их определить?  Нет.  Вам помогали ваши профессиональные навыки?  Наверное, да. Не знаю."
$ws.Cells.Item(9, 2).Value = "  Ммм...  А... Ну, что могло вам помочь?  Ну, классифицирование.  Угу."
$ws.Cells.Item(10, 2).Value = "  Вы любите читать?  Ну, да.  А много читаете?  Когда как. Когда могу зачитываться, когда могу по углам читать."
$ws.Cells.Item(11, 2).Value = "  Могу за месяц схавать 10 книг, а потом такой, ну, мне что-то нечего читать, ну ладно.  Вы решаете задачную логику?  Ну, скорее это, в принципе, как самоцель или как деятельность."
$ws.Cells.Item(12, 2).Value = "  Ну, без разницы.  Ну, скорее всего, да. Потому что приходится.  Угу. Как вы думаете, там могло в"
$ws.Cells.Item(13, 2).Value = "  Не знаю. Ну, в целом.  Есть ли какие-то комментарии к проведению эксперимента?  Ну, что-то не понравилось там."
$ws.Cells.Item(14, 2).Value = "  Да.  Справить что-то можно. Нет?  В принципе, всё понятно.  Хотели бы ещё раз поучаствовать?  Да.  Всё, спасибо."

Write-Host "Applied transcript extension: rows now A1:B14"
